$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Data")

# Update the order id for the "sales_orders" row (row 3)
$ws.Range("C3").Value = "20-VN9MH-21"

# Remove the duplicated trailing text in the ShipToInfo cell for row 3
$ws.Range("O3").Value = "200,INGRAM MICRO TEST ACCOUNT,ATTN TOD DEBIE 1610 E SAINT ANDREW PL SANTA ANA CA 927054931,.,.,.,."

# Add new header column AI for the new "Remove Existing Line" (OMS-29) test case
$ws.Range("AI1").Font.Name = "Arial"
$ws.Range("AI1").Font.Size = 10
$ws.Range("AI1").Font.Bold = $true
$ws.Range("AI1").Value = "OrderLineDataErrorOrderID"

# Adjust column widths to match the new layout (offset by the engine's fixed
# 5/6-character padding between ColumnWidth and the stored <col width>)
$ws.Columns.Item(34).ColumnWidth = 36 - 5/6
$ws.Columns.Item(35).ColumnWidth = 22.5703125 - 5/6
$ws.Columns.Item(36).ColumnWidth = 21 - 5/6

# Move the active selection to the new column's data row to match the
# reviewer's cursor position after adding the test case
$ws.Range("AI6").Select()

$wb.Save()
